# TC01_INS_SpecialTopic-CancerMoonshot.xlsx -- "special topic ins completed"
#
# The workbook tracks 4 DuckDB/SQL queries (StatQuery/TabQuery columns) used
# to populate the INS "Cancer Moonshot" special-topic tab. This change:
#   - finalises the "count" query (adds a trailing semicolon)
#   - fixes the Grants query to pull the grant's own end date
#     (gnt.grant_end_date) instead of re-using the project's end date
#   - fixes the Publications query's title column (pub.publication_title)
#     and adds explicit CASE branches for citation ratios of 1 and 2
#   - moves the active selection to C5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# C2: "Programs/Projects/Grants/Publications" count StatQuery -- finalised
#     with a trailing semicolon.
# ---------------------------------------------------------------------
$countQuery = @"
SELECT DISTINCT
    COUNT(DISTINCT prg.program_id) AS "Programs",
    COUNT(DISTINCT prj.project_id) AS "Projects",
    COUNT(DISTINCT gnt.grant_id) AS "Grants",
    COUNT(DISTINCT pub.pmid) AS "Publications"
FROM 
    df_program prg
LEFT JOIN 
    df_project prj ON prg.program_id = prj."program.program_id"
LEFT JOIN 
    df_grant gnt ON prj.project_id = gnt."project.project_id"
LEFT JOIN 
    df_publication pub ON prj.project_id = pub."project.project_id"
WHERE 
    prg.focus_area LIKE '%Cancer Moonshot%';
"@
$ws.Range("C2").Value = $countQuery

# ---------------------------------------------------------------------
# B4: GrantsTab TabQuery -- now selects the grant's own end date
#     (gnt.grant_end_date) rather than re-using gnt.project_end_date.
# ---------------------------------------------------------------------
$grantQuery = @"
SELECT DISTINCT
    gnt.grant_id AS "Grant ID", 
    prj.project_id AS "Project",
    gnt.grant_title AS "Grant Title",
    gnt.principal_investigators AS "Principal Investigators",
    gnt.program_officers AS "Program Officers",
    gnt.fiscal_year AS "Fiscal Year",
    gnt.grant_end_date AS "Project End Date"
FROM 
    df_grant gnt
LEFT JOIN 
    df_project prj ON gnt."project.project_id" = prj.project_id
LEFT JOIN 
    df_program prg ON prj."program.program_id" = prg.program_id
LEFT JOIN 
    df_publication pub ON prj.project_id = pub."project.project_id"
WHERE 
    prg.focus_area LIKE '%Cancer Moonshot%'
ORDER BY 
    lower(gnt.grant_id) ASC
LIMIT 100;
"@
$ws.Range("B4").Value = $grantQuery

# ---------------------------------------------------------------------
# B5: PublicationsTab TabQuery -- title column renamed to
#     pub.publication_title, and the relative-citation-ratio CASE now
#     special-cases ratios of 1 and 2 (in addition to the existing 0/7).
# ---------------------------------------------------------------------
$pubQuery = @"
SELECT DISTINCT
    pub.pmid AS "PubMed ID", 
    pub.publication_title AS "Title",
    pub.authors AS "Authors",
    pub.publication_date AS "Publication Date",
    pub.cited_by AS "Cited By",
    CASE 
    WHEN pub.relative_citation_ratio = 0 THEN '0'
    WHEN pub.relative_citation_ratio = 7.0 THEN '7'
    WHEN pub.relative_citation_ratio = 1.0 THEN '1'
    WHEN pub.relative_citation_ratio = 2.0 THEN '2'
    WHEN pub.relative_citation_ratio = ROUND(pub.relative_citation_ratio) THEN CAST(ROUND(pub.relative_citation_ratio) AS VARCHAR) 
    ELSE CAST(ROUND(pub.relative_citation_ratio, 2) AS VARCHAR)
END AS "Relative Citation Ratio"
FROM 
    df_publication pub
LEFT JOIN 
    df_project prj ON pub."project.project_id" = prj.project_id
LEFT JOIN 
    df_program prg ON prj."program.program_id" = prg.program_id
LEFT JOIN 
    df_grant gnt ON prj.project_id = gnt."project.project_id"
WHERE 
    prg.focus_area LIKE '%Cancer Moonshot%'
ORDER BY 
    lower(pub.pmid) ASC
LIMIT 100;
"@
$ws.Range("B5").Value = $pubQuery

# ---------------------------------------------------------------------
# Keep the wrap-text / font styling consistent across the TabQuery /
# StatQuery cells (B2:B5, C2) now that they all hold finalised queries.
# ---------------------------------------------------------------------
$queryCells = @("B2", "C2", "B3", "B4", "B5")
foreach ($addr in $queryCells) {
    $rng = $ws.Range($addr)
    $rng.WrapText = $true
    $rng.Font.Size = 12
}

# Special topic work is done -- move off the header row onto the
# finished Publications query cell.
$ws.Range("C5").Select()
